$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add a new column N (2023) mirroring column M's formatting for rows 3-7 ---
# Copy the formats of M3:M7 onto N3:N7 (values for M5/M6 are blank, so this also
# creates the blank-but-styled N5/N6 cells seen in the target).
$ws.Range("M3:M7").Copy()
$ws.Range("N3").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New 2023 column values
$ws.Range("N3").Value = 2023
$ws.Range("N4").Value = 583.5
$ws.Range("N7").Value = 64.08
# N5 and N6 stay blank (same as their M counterparts)

# --- Updated 2022 column (M) values ---
$ws.Range("M4").Value = 923.8
$ws.Range("M7").Value = 64.03

# --- Row height adjustments (Excel auto re-wrapped these rows once the table
#     gained a 14th column) ---
$ws.Rows.Item(1).RowHeight = 29.25
$ws.Rows.Item(4).RowHeight = 26.25
$ws.Rows.Item(5).RowHeight = 25.5
$ws.Rows.Item(6).RowHeight = 28.5
$ws.Rows.Item(7).RowHeight = 41.25

# --- Reset the view so nothing beyond the table is left selected ---
$ws.Range("A1").Select()
